$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,9,13,19,25,31,37,42,48,52,58,64,70,76)

foreach ($r in $rows) {
    $range = $ws.Range("C" + $r + ":T" + $r)
    $range.ClearContents()
}
